$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a temporary text number-format so these numeric-looking strings
# (e.g. "1.010", "0.00000000312") are stored verbatim as text instead of
# being auto-parsed/normalized into numbers by Excel, then restore the
# cell style back to Normal so no stray formatting is left behind.
$changes = @{
    'D2' = '28.552.53'
    'E2' = '  -1.79%  '
    'D3' = '1.964.60'
    'E3' = '  -0.02%  '
    'D4' = '1.010'
    'E4' = '  +0.43%  '
    'D5' = '323.70'
    'E5' = '  -0.97%  '
    'D6' = '1.010'
    'E6' = '  +0.47%  '
    'D7' = '0.4826'
    'E7' = '  -3.55%  '
    'D8' = '0.4082'
    'E8' = '  -3.27%  '
    'D9' = '54.02'
    'E9' = '  +0.36%  '
    'D10' = '0.08527'
    'E10' = '  -6.59%  '
    'E11' = '  -3.12%  '
    'D12' = '22.51'
    'E12' = '  -2.74%  '
    'D13' = '1.954.85'
    'E13' = '  -2.33%  '
    'D14' = '7.624'
    'E14' = '  -3.20%  '
    'D15' = '6.208'
    'E15' = '  -3.53%  '
    'D16' = '1.012'
    'D17' = '91.33'
    'E17' = '  -0.24%  '
    'E18' = '  -2.37%  '
    'D19' = '0.06644'
    'E19' = '  -0.32%  '
    'D20' = '18.58'
    'E20' = '  -3.41%  '
    'E21' = '  +0.46%  '
    'E22' = '  -1.05%  '
    'D23' = '28.614.53'
    'E23' = '  -1.65%  '
    'D24' = '11.54'
    'E24' = '  -3.50%  '
    'D25' = '2.299'
    'E25' = '  +0.68%  '
    'D26' = '2.273.26'
    'E26' = '  +2.24%  '
    'E27' = '  +0.30%  '
    'E28' = '  -1.04%  '
    'D29' = '5.942'
    'E29' = '  -3.29%  '
    'D30' = '2.193'
    'E30' = '  -3.32%  '
    'E31' = '  -1.50%  '
    'D32' = '0.9983'
    'D33' = '0.09732'
    'E33' = '  -1.23%  '
    'D34' = '1.472'
    'E34' = '  -4.30%  '
    'D35' = '5.669'
    'E35' = '  -2.03%  '
    'D36' = '3.694'
    'E36' = '  +0.23%  '
    'D37' = '9.212'
    'E37' = '  +2.99%  '
    'D38' = '0.02346'
    'E38' = '  -3.13%  '
    'D39' = '0.06276'
    'E39' = '  -0.47%  '
    'E40' = '  -2.83%  '
    'D41' = '0.6269'
    'E41' = '  -2.94%  '
    'D42' = '11.31'
    'E42' = '  -1.27%  '
    'E43' = '  +0.49%  '
    'D44' = '0.1927'
    'E44' = '  -3.21%  '
    'D45' = '1.354'
    'E45' = '  +5.44%  '
    'D46' = '13.18'
    'E46' = '  -1.35%  '
    'D47' = '0.5992'
    'E47' = '  -3.55%  '
    'D48' = '2.081'
    'E48' = '  -4.37%  '
    'D49' = '3.417'
    'E49' = '  -1.28%  '
    'E50' = '  -0.90%  '
    'B51' = 'BabyDogeCoin'
    'C51' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D51' = '0.00000000312'
    'E51' = '  -6.27%  '
}

foreach ($ref in $changes.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$ref]
    $cell.Style = "Normal"
}
